$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the project-tracking template table (row 3 header, row 4 sample data)
$ws.Range("A1").Value = "Repo:"
$ws.Range("A3").Value = "Project Name"
$ws.Range("B3").Value = "github Url"
$ws.Range("C3").Value = "Project Language"
$ws.Range("D3").Value = "Time tracking"
$ws.Range("A4").Value = "XX"
$ws.Range("B4").Value = "http://"
$ws.Range("C4").Value = "Java, C++, etc"
$ws.Range("D4").Value = "Yes"

# Fill in the time-tracking URL last
$ws.Range("B1").Value = "https://issues.apache.org/jira/secure/Dashboard.jspa"

# Turn the sample github URL into a real hyperlink (applies the built-in Hyperlink style)
$ws.Hyperlinks.Add($ws.Range("B4"), "http://") | Out-Null

# Leave the selection where the author last clicked
$ws.Range("C7").Select() | Out-Null
